$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "1024_3_000033_h"
$ws.Range("C2").Value = 82
$ws.Range("D2").Value = 0.8661446041531033
$ws.Range("E2").Value = 0.8445808004347161
$ws.Range("F2").Value = 0.9243334039168646
$ws.Range("G2").Value = 0.2034571644607057
$ws.Range("H2").Value = 0.8826592541663918
$ws.Range("I2").Value = 0.9323465529848501
